$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C (coin name / link) are plain text - safe to assign directly.
# Columns D and E (price / volume%) look numeric, so the sheet stores them as literal
# text (e.g. "266.13", "2.02%") preserving exact formatting/precision. Force the cell
# to Text format first so Excel does not silently convert them to numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "266.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.02%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.84%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.695"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.49%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06120"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.93%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.741"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.08%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8500"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9092"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.98%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1407"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.40%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04812"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.81%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07086"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.12%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03125"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.56%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09041"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.03%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001542"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.49%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006192"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.64%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006036"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.82%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.451"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.05%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.162"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.21%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.147"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.47%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3072"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.19%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.55%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.116"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.58%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04237"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.44%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001180"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.00%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004070"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "7.04%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001198"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.18%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03940"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.65%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1120"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.81%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004169"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.68%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002106"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.75%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01171"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-28.26%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005076"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.73%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.18%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2516"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "56.21%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002096"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.18%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001996"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.18%"
